$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing content/row-heights so we can rebuild the sheet cleanly ---
$ws.Range("A1:C25").Clear()
for ($r = 1; $r -le 25; $r++) {
  $ws.Rows.Item($r).RowHeight = 15
  $ws.Rows.Item($r).AutoFit()
}

# --- Write cell values row by row ---
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOM3259'
$ws.Range("C2").Value = 'LOM3259'

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Materiais e Dispositivos Eletrônicos'
$ws.Range("C3").Value = ' Materiais e Dispositivos Eletrônicos'

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Electronic Materials and Devices'
$ws.Range("C4").Value = 'Electronic Materials and Devices'

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2023'
$ws.Range("C8").Value = '01/01/2023'

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EF-8'
$ws.Range("C9").Value = 'EF-8'

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Propiciar ao aluno os conhecimentos básicos de materiais eletrônicos visando sua aplicação em dispositivos.'
$ws.Range("C10").Value = 'Propiciar ao aluno os conhecimentos básicos de materiais eletrônicos visando sua aplicação em dispositivos.'

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'Provide the student with the basic knowledge of electronic materials aiming their application in devices.'
$ws.Range("C11").Value = 'Provide the student with the basic knowledge of electronic materials aiming their application in devices.'

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("B13").Value = '144651 - Antonio Fernando Sartori'
$ws.Range("C13").Value = '144651 - Antonio Fernando Sartori'

$ws.Range("B14").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C14").Value = '5840730 - Antonio Jefferson da Silva Machado'

$ws.Range("B15").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C15").Value = '519033 - Carlos Yujiro Shigue'

$ws.Range("A16").Value = 'Programa resumido:'
$ws.Range("B16").Value = 'Materiais para eletrônica. Eletrônica e Física do Estado Sólido. Materiais e dispositivos semicondutores. Materiais e dispositivos optoeletrônicos. Materiais e dispositivos dielétricos e piezelétricos.'
$ws.Range("C16").Value = 'Materiais para eletrônica. Eletrônica e Física do Estado Sólido. Materiais e dispositivos semicondutores. Materiais e dispositivos optoeletrônicos. Materiais e dispositivos dielétricos e piezelétricos.'

$ws.Range("A17").Value = 'Short syllabus:'
$ws.Range("B17").Value = 'Materials for electronics. Electronics and Solid State Physics. Semiconductor materials and devices. Optoelectronic materials and devices. Dielectric and piezoelectric materials and devices.'
$ws.Range("C17").Value = 'Materials for electronics. Electronics and Solid State Physics. Semiconductor materials and devices. Optoelectronic materials and devices. Dielectric and piezoelectric materials and devices.'

$ws.Range("A18").Value = 'Programa:'
$ws.Range("B18").Value = 'Materiais para aplicações eletrônicas: metais, cerâmicas, vidros e polímeros. Monocristais e filmes finos.Ondas e partículas na matéria. Elétrons em átomos e cristais. Estruturas de bandas de energia. Propriedades eletrônicas e espectroscópicas de materiais. Materiais condutores, semicondutores e isolantes. Propriedades eletrônicas em semicondutores. Transporte elétrico. Dispositivos semicondutores. Junção pn. Contato metal-semicondutor e semicondutor-isolante. Dispositivos semicondutores: diodos e transistores bipolares e FET. Materiais e dispositivos optoeletrônicos. LED, laser semicondutor, fotodetetores e células fotovoltaicas. Tipos e propriedades dos materiais dielétricos. Materiais ferroelétricos e piezelétricos. Dispositivos baseados em materiais dielétricos e piezelétricos. Aplicações.'
$ws.Range("C18").Value = 'Materiais para aplicações eletrônicas: metais, cerâmicas, vidros e polímeros. Monocristais e filmes finos.Ondas e partículas na matéria. Elétrons em átomos e cristais. Estruturas de bandas de energia. Propriedades eletrônicas e espectroscópicas de materiais. Materiais condutores, semicondutores e isolantes. Propriedades eletrônicas em semicondutores. Transporte elétrico. Dispositivos semicondutores. Junção pn. Contato metal-semicondutor e semicondutor-isolante. Dispositivos semicondutores: diodos e transistores bipolares e FET. Materiais e dispositivos optoeletrônicos. LED, laser semicondutor, fotodetetores e células fotovoltaicas. Tipos e propriedades dos materiais dielétricos. Materiais ferroelétricos e piezelétricos. Dispositivos baseados em materiais dielétricos e piezelétricos. Aplicações.'

$ws.Range("A19").Value = 'Syllabus:'
$ws.Range("B19").Value = 'Materials for electronic applications: metals, ceramics, glasses and polymers. Single crystals and thin films.Waves and particles in matter. Electrons in atoms and crystals. Energy band structures. Electronic and spectroscopic properties of materials.Conducting, semiconducting and insulating materials. Electronic properties in semiconductors. Electric transport. Semiconductor devices. pn junction Metal-semiconductor and semiconductor-insulator contact. Semiconductor devices: diodes and bipolar and FET transistors.Optoelectronic materials and devices. LED, semiconductor laser, photodetectors and photovoltaic cells.Types and properties of dielectric materials. Ferroelectric and piezoelectric materials. Devices based on dielectric and piezoelectric materials. Applications.'
$ws.Range("C19").Value = 'Materials for electronic applications: metals, ceramics, glasses and polymers. Single crystals and thin films.Waves and particles in matter. Electrons in atoms and crystals. Energy band structures. Electronic and spectroscopic properties of materials.Conducting, semiconducting and insulating materials. Electronic properties in semiconductors. Electric transport. Semiconductor devices. pn junction Metal-semiconductor and semiconductor-insulator contact. Semiconductor devices: diodes and bipolar and FET transistors.Optoelectronic materials and devices. LED, semiconductor laser, photodetectors and photovoltaic cells.Types and properties of dielectric materials. Ferroelectric and piezoelectric materials. Devices based on dielectric and piezoelectric materials. Applications.'

$ws.Range("A20").Value = 'Avaliação:'

$ws.Range("A21").Value = 'Método:'
$ws.Range("B21").Value = 'Aulas expositivas, práticas, seminários e exercícios.'
$ws.Range("C21").Value = 'Aulas expositivas, práticas, seminários e exercícios.'

$ws.Range("A22").Value = 'Critério:'
$ws.Range("B22").Value = 'Média das notas de provas, relatórios e apresentações.'
$ws.Range("C22").Value = 'Média das notas de provas, relatórios e apresentações.'

$ws.Range("A23").Value = 'Norma de recuperação:'
$ws.Range("B23").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C23").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

$ws.Range("A24").Value = 'Bibliografia:'
$ws.Range("B24").Value = 'REZENDE, S. M. Materiais e Dispositivos Eletrônicos, São Paulo: Livraria da Física, 2014.
SWART, J. W. Semicondutores - Fundamentos, Técnicas e Aplicações, Campinas: Editora da Unicamp, 2008.
YU, P. Y.; CARDONA, M. Fundamentals of Semiconductors: Physics and Materials Properties, Springer, 2005.
KWOK, H. L. Electronic Materials, Boston: PWS Publishing, 1997.
HORENSTEIN, M. N. Microeletrônica: Circuitos & Dispositivos. Rio de Janeiro, Prentice Hall do Brasil, 1996.
SCHMIDT, W. Materiais Elétricos, vol. I, Ed. Edgard Blücher, SP, 1998.
SCHMIDT, W. Materiais Elétricos, vol. II, Ed. Edgard Blücher, SP, 1995.
HIPPEL, A. R. Dielectric Materials and Applications, Artech House, 1995.
CHOUDHARY, R. N. Dielectric Materials: Introduction, Research and Applications, Nova Science Pub., 2009.
YANG, J. An Introduction to Theory of Piezoelectricity, Springer, 2004. 
VIVES, A. A. Piezoelectric Transducer and Applications, Springer, 2008.'
$ws.Range("C24").Value = 'REZENDE, S. M. Materiais e Dispositivos Eletrônicos, São Paulo: Livraria da Física, 2014.
SWART, J. W. Semicondutores - Fundamentos, Técnicas e Aplicações, Campinas: Editora da Unicamp, 2008.
YU, P. Y.; CARDONA, M. Fundamentals of Semiconductors: Physics and Materials Properties, Springer, 2005.
KWOK, H. L. Electronic Materials, Boston: PWS Publishing, 1997.
HORENSTEIN, M. N. Microeletrônica: Circuitos & Dispositivos. Rio de Janeiro, Prentice Hall do Brasil, 1996.
SCHMIDT, W. Materiais Elétricos, vol. I, Ed. Edgard Blücher, SP, 1998.
SCHMIDT, W. Materiais Elétricos, vol. II, Ed. Edgard Blücher, SP, 1995.
HIPPEL, A. R. Dielectric Materials and Applications, Artech House, 1995.
CHOUDHARY, R. N. Dielectric Materials: Introduction, Research and Applications, Nova Science Pub., 2009.
YANG, J. An Introduction to Theory of Piezoelectricity, Springer, 2004. 
VIVES, A. A. Piezoelectric Transducer and Applications, Springer, 2008.'

$ws.Range("A25").Value = 'Requisitos:'

$ws.Range("B26").Value = 'LOM3215 -  Física do Estado Sólido  (Requisito)
'
$ws.Range("C26").Value = 'LOM3215 -  Física do Estado Sólido  (Requisito)
'

$ws.Range("B27").Value = 'LOM3234 -  Óptica Física  (Requisito)
'
$ws.Range("C27").Value = 'LOM3234 -  Óptica Física  (Requisito)
'

$ws.Range("B28").Value = 'LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)
'
$ws.Range("C28").Value = 'LOM3263 -  Eletrônica Fundamental e Aplicada  (Requisito)
'

# --- Apply custom row heights ---
$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).RowHeight = 120
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 60
$ws.Rows.Item(24).RowHeight = 120
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30

# --- Fix column width / cols split (A vs B) ---
$ws.Columns.Item(1).ColumnWidth = 29.84

